$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-17T07:09:40.579586"
    3 = "2025-10-17T07:09:40.579586"
    4 = "2025-10-17T07:09:40.579586"
    5 = "2025-10-17T07:09:40.579586"
    6 = "2025-10-17T07:09:40.579586"
    7 = "2025-10-17T07:09:40.579586"
    8 = "2025-10-17T07:09:40.579586"
    9 = "2025-10-17T07:09:40.579586"
    10 = "2025-10-17T07:09:40.579586"
    11 = "2025-10-17T07:09:40.579586"
    12 = "2025-10-17T07:09:40.579586"
    13 = "2025-10-17T07:09:40.579586"
    14 = "2025-10-17T07:09:40.579586"
    15 = "2025-10-17T07:09:40.579586"
    16 = "2025-10-17T07:09:40.579586"
    17 = "2025-10-17T07:09:40.579586"
    18 = "2025-10-17T07:09:40.579586"
    19 = "2025-10-17T07:09:40.579586"
    20 = "2025-10-17T07:09:40.579586"
    21 = "2025-10-17T07:09:40.579586"
    22 = "2025-10-17T07:09:40.579586"
    23 = "2025-10-17T07:09:40.579586"
    24 = "2025-10-17T07:09:40.579586"
    25 = "2025-10-17T07:09:40.579586"
    26 = "2025-10-17T07:09:40.579586"
    27 = "2025-10-17T07:09:40.579586"
    28 = "2025-10-17T07:09:40.579586"
    29 = "2025-10-17T07:09:40.579586"
    30 = "2025-10-17T07:09:40.579586"
    31 = "2025-10-17T07:09:40.579586"
    32 = "2025-10-17T07:09:40.579586"
    33 = "2025-10-17T07:09:40.579586"
    34 = "2025-10-17T07:09:40.579586"
    35 = "2025-10-17T07:09:40.579586"
    36 = "2025-10-17T07:09:40.579586"
    37 = "2025-10-17T07:09:40.579586"
    38 = "2025-10-17T07:09:40.579586"
    39 = "2025-10-17T07:09:40.579586"
    40 = "2025-10-17T07:09:40.579586"
    41 = "2025-10-17T07:09:40.579586"
    42 = "2025-10-17T07:09:40.590593"
    43 = "2025-10-17T07:09:40.590593"
    44 = "2025-10-17T07:09:40.590593"
    45 = "2025-10-17T07:09:40.590593"
    46 = "2025-10-17T07:09:40.690424"
    47 = "2025-10-17T07:09:40.690424"
    48 = "2025-10-17T07:09:40.690424"
    49 = "2025-10-17T07:09:40.706053"
    50 = "2025-10-17T07:09:40.706199"
    51 = "2025-10-17T07:09:40.706199"
    52 = "2025-10-17T07:09:40.706199"
    53 = "2025-10-17T07:09:40.706199"
    54 = "2025-10-17T07:09:40.707197"
    55 = "2025-10-17T07:09:40.707197"
    56 = "2025-10-17T07:09:40.707197"
    57 = "2025-10-17T07:09:40.707197"
    58 = "2025-10-17T07:09:40.707197"
    59 = "2025-10-17T07:09:40.707197"
    60 = "2025-10-17T07:09:40.707197"
    61 = "2025-10-17T07:09:40.708199"
    62 = "2025-10-17T07:09:40.708199"
    63 = "2025-10-17T07:09:40.708199"
    64 = "2025-10-17T07:09:40.708199"
    65 = "2025-10-17T07:09:40.708199"
    66 = "2025-10-17T07:09:40.708199"
    67 = "2025-10-17T07:09:40.708199"
    68 = "2025-10-17T07:09:40.708199"
    69 = "2025-10-17T07:09:40.708199"
    70 = "2025-10-17T07:09:40.708199"
    71 = "2025-10-17T07:09:40.709196"
    72 = "2025-10-17T07:09:40.709196"
    73 = "2025-10-17T07:09:40.709196"
    74 = "2025-10-17T07:09:40.709196"
    75 = "2025-10-17T07:09:40.791085"
    76 = "2025-10-17T07:09:40.791085"
    77 = "2025-10-17T07:09:40.791085"
    78 = "2025-10-17T07:09:40.791085"
    79 = "2025-10-17T07:09:40.791085"
    80 = "2025-10-17T07:09:40.791085"
    81 = "2025-10-17T07:09:40.791085"
    82 = "2025-10-17T07:09:40.791085"
    83 = "2025-10-17T07:09:40.791085"
    84 = "2025-10-17T07:09:40.791085"
    85 = "2025-10-17T07:09:40.791085"
    86 = "2025-10-17T07:09:40.791085"
    87 = "2025-10-17T07:09:40.791085"
    88 = "2025-10-17T07:09:40.791085"
    89 = "2025-10-17T07:09:40.791085"
    90 = "2025-10-17T07:09:40.791085"
    91 = "2025-10-17T07:09:40.791085"
    92 = "2025-10-17T07:09:40.791085"
    93 = "2025-10-17T07:09:40.791085"
    94 = "2025-10-17T07:09:40.791085"
    95 = "2025-10-17T07:09:40.791085"
    96 = "2025-10-17T07:09:40.791085"
    97 = "2025-10-17T07:09:40.791085"
    98 = "2025-10-17T07:09:40.791085"
    99 = "2025-10-17T07:09:40.791085"
    100 = "2025-10-17T07:09:40.791085"
    101 = "2025-10-17T07:09:40.791085"
    102 = "2025-10-17T07:09:40.791085"
    103 = "2025-10-17T07:09:40.879479"
    104 = "2025-10-17T07:09:40.879479"
    105 = "2025-10-17T07:09:40.879479"
    106 = "2025-10-17T07:09:40.879479"
    107 = "2025-10-17T07:09:40.879479"
    108 = "2025-10-17T07:09:40.879479"
    109 = "2025-10-17T07:09:40.879479"
    110 = "2025-10-17T07:09:40.879479"
    111 = "2025-10-17T07:09:40.879479"
    112 = "2025-10-17T07:09:40.879479"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item([int]$row, 26).Value = $timestamps[$row]
}
